# Apply crypto price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.935.51'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '3.037.29'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.12'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.90'
$ws.Range('E6').Value = '  +7.18%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.032.91'
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('E10').Value = '  +13.54%  '
$ws.Range('E11').Value = '  +4.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.465'
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.88'
$ws.Range('E14').Value = '  +4.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '3.541.43'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').Value = '62.935.75'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').Value = '3.037.99'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '454.60'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.31'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.699'
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('E23').Value = '  +3.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.11'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.32'
$ws.Range('E25').Value = '  +7.28%  '
$ws.Range('E26').Value = '  +4.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.43'
$ws.Range('E27').Value = '  +4.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('E29').Value = '  +5.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.26'
$ws.Range('E30').Value = '  +10.45%  '
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.66'
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('D35').Value = '0.0₃0858'
$ws.Range('E35').Value = '  +4.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  +2.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.95'
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.19'
$ws.Range('E38').Value = '  +10.28%  '
$ws.Range('E39').Value = '  +7.70%  '
$ws.Range('E40').Value = '  +2.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.37'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.306'
$ws.Range('E43').Value = '  +13.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.92'
$ws.Range('E44').Value = '  +10.59%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '392.82'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0363'
$ws.Range('E46').Value = '  +3.57%  '
$ws.Range('D47').Value = '2.723.42'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.84'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  +7.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.58'
$ws.Range('E51').Value = '  +4.91%  '
